$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 106999.71
$ws.Range("I100").Value = 125787
$ws.Range("J100").Value = 81950
$ws.Range("K100").Value = 125787
$ws.Range("L100").Value = 81950
$ws.Range("M100").Value = -125246
$ws.Range("N100").Value = -83032
$ws.Range("H113").Value = 8949.25
$ws.Range("I113").Value = 6936.5
$ws.Range("K113").Value = 6936.5
$ws.Range("M113").Value = -3682.5
$ws.Range("H135").Value = 4872.3076
$ws.Range("J135").Value = 3696.4736
$ws.Range("L135").Value = 33268.2624
$ws.Range("N135").Value = -38338.2624
$ws.Range("H137").Value = 8024.9395
$ws.Range("I137").Value = 12792.611
$ws.Range("J137").Value = 2303.7334
$ws.Range("K137").Value = 38377.833
$ws.Range("L137").Value = 6911.2002
$ws.Range("M137").Value = -35827.833
$ws.Range("N137").Value = -12011.2002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 88477.625
$ws.Range("I45").Value = 137977.53
$ws.Range("J45").Value = 5977.778
$ws.Range("K45").Value = 137977.53
$ws.Range("L45").Value = 5977.778
$ws.Range("M45").Value = -137600.53
$ws.Range("N45").Value = -6731.778
$ws.Range("H61").Value = 8791.5
$ws.Range("I61").Value = 10643.737
$ws.Range("J61").Value = 4881.222
$ws.Range("K61").Value = 10643.737
$ws.Range("L61").Value = 4881.222
$ws.Range("M61").Value = -10431.737
$ws.Range("N61").Value = -5305.222
$ws.Range("H63").Value = 2548.5454
$ws.Range("J63").Value = 2998
$ws.Range("L63").Value = 2998
$ws.Range("N63").Value = -4370
$ws.Range("H66").Value = 2548.5454
$ws.Range("J66").Value = 2998
$ws.Range("L66").Value = 14990
$ws.Range("N66").Value = -21854
$ws.Range("H110").Value = 2146.625
$ws.Range("I110").Value = 1296.2354
$ws.Range("J110").Value = 4211.857
$ws.Range("K110").Value = 1296.2354
$ws.Range("L110").Value = 4211.857
$ws.Range("M110").Value = 748.7646
$ws.Range("N110").Value = -8301.857
$ws.Range("H122").Value = 2507558.8
$ws.Range("I122").Value = 7588.375
$ws.Range("K122").Value = 22765.125
$ws.Range("M122").Value = -20315.125
$ws.Range("H132").Value = 4401.6523
$ws.Range("I132").Value = 4230.9414
$ws.Range("J132").Value = 4885.3335
$ws.Range("K132").Value = 12692.8242
$ws.Range("L132").Value = 14656.0005
$ws.Range("M132").Value = -10162.8242
$ws.Range("N132").Value = -19716.0005
$ws.Range("H136").Value = 8791.5
$ws.Range("I136").Value = 10643.737
$ws.Range("J136").Value = 4881.222
$ws.Range("K136").Value = 31931.211
$ws.Range("L136").Value = 14643.666
$ws.Range("M136").Value = -29381.211
$ws.Range("N136").Value = -19743.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 10000
$ws.Range("I44").Value = 10000
$ws.Range("K44").Value = 10000
$ws.Range("M44").Value = -9503
$ws.Range("H64").Value = 11063.375
$ws.Range("I64").Value = 16000
$ws.Range("K64").Value = 16000
$ws.Range("M64").Value = -15775
$ws.Range("H67").Value = 11063.375
$ws.Range("I67").Value = 16000
$ws.Range("K67").Value = 16000
$ws.Range("M67").Value = -15220
$ws.Range("H94").Value = 8023.3
$ws.Range("I94").Value = 10315.607
$ws.Range("K94").Value = 10315.607
$ws.Range("M94").Value = -9864.607
$ws.Range("H134").Value = 14311.2
$ws.Range("I134").Value = 15590.223
$ws.Range("K134").Value = 46770.669
$ws.Range("M134").Value = -44235.669

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 17087.334
$ws.Range("I7").Value = 40330
$ws.Range("J7").Value = 485.42856
$ws.Range("K7").Value = 40330
$ws.Range("L7").Value = 485.42856
$ws.Range("M7").Value = -40217
$ws.Range("N7").Value = -711.4285600000001
$ws.Range("H31").Value = 11180.066
$ws.Range("I31").Value = 12700.272
$ws.Range("J31").Value = 6999.5
$ws.Range("K31").Value = 12700.272
$ws.Range("L31").Value = 6999.5
$ws.Range("M31").Value = -12405.272
$ws.Range("N31").Value = -7589.5
$ws.Range("H34").Value = 11180.066
$ws.Range("I34").Value = 12700.272
$ws.Range("J34").Value = 6999.5
$ws.Range("K34").Value = 12700.272
$ws.Range("L34").Value = 6999.5
$ws.Range("M34").Value = -12498.272
$ws.Range("N34").Value = -7403.5
$ws.Range("H58").Value = 3248.074
$ws.Range("I58").Value = 3519.0625
$ws.Range("J58").Value = 2853.9092
$ws.Range("K58").Value = 3519.0625
$ws.Range("L58").Value = 2853.9092
$ws.Range("M58").Value = -3316.0625
$ws.Range("N58").Value = -3259.9092
$ws.Range("H81").Value = 48000
$ws.Range("I81").Value = 48000
$ws.Range("K81").Value = 48000
$ws.Range("M81").Value = -47002
$ws.Range("H84").Value = 48000
$ws.Range("I84").Value = 48000
$ws.Range("K84").Value = 144000
$ws.Range("M84").Value = -139008
$ws.Range("H94").Value = 2311.75
$ws.Range("J94").Value = 2311.75
$ws.Range("L94").Value = 2311.75
$ws.Range("N94").Value = -3213.75
$ws.Range("H99").Value = 11617311
$ws.Range("I99").Value = 38706036
$ws.Range("J99").Value = 7857.143
$ws.Range("K99").Value = 38706036
$ws.Range("L99").Value = 7857.143
$ws.Range("M99").Value = -38704538
$ws.Range("N99").Value = -10853.143
$ws.Range("H126").Value = 11617311
$ws.Range("I126").Value = 38706036
$ws.Range("J126").Value = 7857.143
$ws.Range("K126").Value = 116118108
$ws.Range("L126").Value = 23571.429
$ws.Range("M126").Value = -116115638
$ws.Range("N126").Value = -28511.429
$ws.Range("H132").Value = 1308.1143
$ws.Range("I132").Value = 1255.7812
$ws.Range("K132").Value = 3767.3436
$ws.Range("M132").Value = -1237.3436
$ws.Range("H136").Value = 3248.074
$ws.Range("I136").Value = 3519.0625
$ws.Range("J136").Value = 2853.9092
$ws.Range("K136").Value = 10557.1875
$ws.Range("L136").Value = 8561.7276
$ws.Range("M136").Value = -8007.1875
$ws.Range("N136").Value = -13661.7276
$ws.Range("H137").Value = 37500
$ws.Range("I137").Value = 25000
$ws.Range("K137").Value = 25000
$ws.Range("M137").Value = -19900

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 4500
$ws.Range("J69").Value = 4500
$ws.Range("L69").Value = 13500
$ws.Range("N69").Value = -15122
$ws.Range("H72").Value = 4500
$ws.Range("J72").Value = 4500
$ws.Range("L72").Value = 40500
$ws.Range("N72").Value = -48612
$ws.Range("H80").Value = 86756.2
$ws.Range("I80").Value = 3785
$ws.Range("J80").Value = 95975.22
$ws.Range("K80").Value = 11355
$ws.Range("L80").Value = 287925.66
$ws.Range("M80").Value = -10419
$ws.Range("N80").Value = -289797.66
$ws.Range("H83").Value = 86756.2
$ws.Range("I83").Value = 3785
$ws.Range("J83").Value = 95975.22
$ws.Range("K83").Value = 34065
$ws.Range("L83").Value = 863776.98
$ws.Range("M83").Value = -29385
$ws.Range("N83").Value = -873136.98
$ws.Range("H93").Value = 11499.667
$ws.Range("J93").Value = 13399.8
$ws.Range("L93").Value = 40199.39999999999
$ws.Range("N93").Value = -43943.39999999999
$ws.Range("H107").Value = 2192.889
$ws.Range("J107").Value = 2192.889
$ws.Range("L107").Value = 6578.667
$ws.Range("N107").Value = -10418.667
$ws.Range("H113").Value = 11983.333
$ws.Range("I113").Value = 916.6667
$ws.Range("K113").Value = 2750.0001
$ws.Range("M113").Value = -580.0001000000002
$ws.Range("H122").Value = 4895.3784
$ws.Range("J122").Value = 5287.8486
$ws.Range("L122").Value = 47590.6374
$ws.Range("N122").Value = -52490.6374
$ws.Range("H131").Value = 11366429
$ws.Range("I131").Value = 250019800
$ws.Range("J131").Value = 1982
$ws.Range("K131").Value = 750059400
$ws.Range("L131").Value = 5946
$ws.Range("M131").Value = -750054360
$ws.Range("N131").Value = -16026
$ws.Range("H132").Value = 30877.941
$ws.Range("I132").Value = 746
$ws.Range("K132").Value = 6714
$ws.Range("M132").Value = -4184
$ws.Range("H137").Value = 4349.6924
$ws.Range("I137").Value = 1584.8
$ws.Range("K137").Value = 4754.4
$ws.Range("M137").Value = 345.6000000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 7895.5
$ws.Range("I102").Value = 9228.588
$ws.Range("K102").Value = 9228.588
$ws.Range("M102").Value = -7606.588
$ws.Range("H113").Value = 9848.786
$ws.Range("I113").Value = 14999.5
$ws.Range("J113").Value = 2981.1667
$ws.Range("K113").Value = 14999.5
$ws.Range("L113").Value = 2981.1667
$ws.Range("M113").Value = -12829.5
$ws.Range("N113").Value = -7321.1667
$ws.Range("H132").Value = 5766.375
$ws.Range("I132").Value = 6228.3335
$ws.Range("J132").Value = 2532.6667
$ws.Range("K132").Value = 18685.0005
$ws.Range("L132").Value = 7598.000100000001
$ws.Range("M132").Value = -16155.0005
$ws.Range("N132").Value = -12658.0001
$ws.Range("H136").Value = 27395.5
$ws.Range("J136").Value = 27395.5
$ws.Range("L136").Value = 82186.5
$ws.Range("N136").Value = -87286.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 19635.5
$ws.Range("I40").Value = 22642.863
$ws.Range("K40").Value = 22642.863
$ws.Range("M40").Value = -22506.863
$ws.Range("H93").Value = 4151.4736
$ws.Range("I93").Value = 4405.1763
$ws.Range("K93").Value = 4405.1763
$ws.Range("M93").Value = -3157.1763
$ws.Range("H132").Value = 1492699.1
$ws.Range("I132").Value = 2484299.2
$ws.Range("K132").Value = 7452897.600000001
$ws.Range("M132").Value = -7450367.600000001
$ws.Range("H136").Value = 5211.4546
$ws.Range("I136").Value = 4447
$ws.Range("K136").Value = 13341
$ws.Range("M136").Value = -10791

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 34000000
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 34000000
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 34000000
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -34000616
$ws.Range("H136").Value = 412889.53
$ws.Range("I136").Value = 502225.88
$ws.Range("K136").Value = 1506677.64
$ws.Range("H139").Value = 92600
$ws.Range("J139").Value = 63250
$ws.Range("L139").Value = 63250
$ws.Range("N139").Value = -73530
